$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before D, shifting existing D:K data to F:M
$ws.Range("D:E").Insert()

# Copy number formatting from the (now-shifted) old D,E columns (F,G) into the new D,E columns
$ws.Range("F5:F102").Copy()
$ws.Range("D5:D102").PasteSpecial(-4122)
$ws.Range("G5:G102").Copy()
$ws.Range("E5:E102").PasteSpecial(-4122)

# Populate values for all affected cells (new D,E columns plus any restated figures in F..M)
$ws.Range("D7").Value = 43463
$ws.Range("E7").Value = 43372
$ws.Range("F7").Value = 43281
$ws.Range("G7").Value = 43190
$ws.Range("H7").Value = 43099
$ws.Range("I7").Value = 43008
$ws.Range("J7").Value = 42917
$ws.Range("K7").Value = 42826
$ws.Range("L7").Value = 42735
$ws.Range("M7").Value = 42643
$ws.Range("D8").Value = 115300
$ws.Range("E8").Value = 104200
$ws.Range("F8").Value = 104000
$ws.Range("G8").Value = 94500
$ws.Range("H8").Value = 94000
$ws.Range("I8").Value = 91300
$ws.Range("J8").Value = 83700
$ws.Range("K8").Value = 64100
$ws.Range("L8").Value = 63000
$ws.Range("M8").Value = 60100
$ws.Range("D9").Value = 61200
$ws.Range("E9").Value = 54200
$ws.Range("F9").Value = 53300
$ws.Range("G9").Value = 47900
$ws.Range("H9").Value = 48600
$ws.Range("I9").Value = 44500
$ws.Range("J9").Value = 42400
$ws.Range("K9").Value = 32300
$ws.Range("L9").Value = 32100
$ws.Range("M9").Value = 32400
$ws.Range("D10").Value = 54100
$ws.Range("E10").Value = 50000
$ws.Range("F10").Value = 50700
$ws.Range("G10").Value = 46600
$ws.Range("H10").Value = 45500
$ws.Range("I10").Value = 46800
$ws.Range("J10").Value = 41300
$ws.Range("K10").Value = 31800
$ws.Range("L10").Value = 30900
$ws.Range("M10").Value = 27700
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("F12").Value = "NA"
$ws.Range("G12").Value = "NA"
$ws.Range("H12").Value = "NA"
$ws.Range("I12").Value = "NA"
$ws.Range("J12").Value = "NA"
$ws.Range("K12").Value = "NA"
$ws.Range("L12").Value = "NA"
$ws.Range("M12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("F13").Value = 0
$ws.Range("G13").Value = 0
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 0
$ws.Range("F14").Value = 0
$ws.Range("G14").Value = 0
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = 0
$ws.Range("D15").Value = 5700
$ws.Range("E15").Value = 4100
$ws.Range("F15").Value = 3800
$ws.Range("G15").Value = 3800
$ws.Range("H15").Value = 3600
$ws.Range("I15").Value = 3800
$ws.Range("J15").Value = 3300
$ws.Range("K15").Value = 2500
$ws.Range("L15").Value = 1900
$ws.Range("M15").Value = 1600
$ws.Range("D17").Value = 105300
$ws.Range("E17").Value = 94200
$ws.Range("F17").Value = 94700
$ws.Range("G17").Value = 88200
$ws.Range("H17").Value = 85800
$ws.Range("I17").Value = 82300
$ws.Range("J17").Value = 76900
$ws.Range("K17").Value = 61500
$ws.Range("L17").Value = 58000
$ws.Range("M17").Value = 54600
$ws.Range("D18").Value = 10000
$ws.Range("E18").Value = 10000
$ws.Range("F18").Value = 9300
$ws.Range("G18").Value = 6300
$ws.Range("H18").Value = 8200
$ws.Range("I18").Value = 9000
$ws.Range("J18").Value = 6800
$ws.Range("K18").Value = 2600
$ws.Range("L18").Value = 5000
$ws.Range("M18").Value = 5500
$ws.Range("D20").Value = 0
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 0
$ws.Range("G20").Value = 0
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = 0
$ws.Range("D21").Value = 15800
$ws.Range("E21").Value = 14000
$ws.Range("F21").Value = 13200
$ws.Range("G21").Value = 10100
$ws.Range("H21").Value = 11800
$ws.Range("I21").Value = 12700
$ws.Range("J21").Value = 10100
$ws.Range("K21").Value = 5000
$ws.Range("L21").Value = 7000
$ws.Range("M21").Value = 7100
$ws.Range("D22").Value = 300
$ws.Range("E22").Value = 500
$ws.Range("F22").Value = 700
$ws.Range("G22").Value = 600
$ws.Range("H22").Value = 900
$ws.Range("I22").Value = 500
$ws.Range("J22").Value = 300
$ws.Range("K22").Value = 200
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = 100
$ws.Range("D23").Value = 9800
$ws.Range("E23").Value = 9500
$ws.Range("F23").Value = 8700
$ws.Range("G23").Value = 5700
$ws.Range("H23").Value = 7300
$ws.Range("I23").Value = 8400
$ws.Range("J23").Value = 6600
$ws.Range("K23").Value = 2300
$ws.Range("L23").Value = 5000
$ws.Range("M23").Value = 5400
$ws.Range("D24").Value = 2000
$ws.Range("E24").Value = 2200
$ws.Range("F24").Value = 1100
$ws.Range("G24").Value = 1400
$ws.Range("H24").Value = 1700
$ws.Range("I24").Value = 2500
$ws.Range("J24").Value = 2300
$ws.Range("K24").Value = 0
$ws.Range("L24").Value = 1700
$ws.Range("M24").Value = 2000
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 0
$ws.Range("G25").Value = 0
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = 0
$ws.Range("D26").Value = 7800
$ws.Range("E26").Value = 7300
$ws.Range("F26").Value = 7600
$ws.Range("G26").Value = 4300
$ws.Range("H26").Value = 5600
$ws.Range("I26").Value = 5900
$ws.Range("J26").Value = 4300
$ws.Range("K26").Value = 2300
$ws.Range("L26").Value = 3300
$ws.Range("M26").Value = 3400
$ws.Range("D27").Value = 7800
$ws.Range("E27").Value = 7300
$ws.Range("F27").Value = 7600
$ws.Range("G27").Value = 4300
$ws.Range("H27").Value = 5600
$ws.Range("I27").Value = 5900
$ws.Range("J27").Value = 4300
$ws.Range("K27").Value = 2300
$ws.Range("L27").Value = 3300
$ws.Range("M27").Value = 3400
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("F28").Value = 0
$ws.Range("G28").Value = 0
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("M28").Value = 0
$ws.Range("D29").Value = -100
$ws.Range("E29").Value = "NA"
$ws.Range("F29").Value = "NA"
$ws.Range("G29").Value = "NA"
$ws.Range("H29").Value = 5900
$ws.Range("I29").Value = "NA"
$ws.Range("J29").Value = "NA"
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = 0
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("F30").Value = 0
$ws.Range("G30").Value = 0
$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = 0
$ws.Range("M30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("F31").Value = 0
$ws.Range("G31").Value = 0
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = 0
$ws.Range("D32").Value = 0
$ws.Range("E32").Value = 0
$ws.Range("F32").Value = 0
$ws.Range("G32").Value = 0
$ws.Range("H32").Value = 0
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = 0
$ws.Range("D33").Value = 7700
$ws.Range("E33").Value = 7300
$ws.Range("F33").Value = 7600
$ws.Range("G33").Value = 4300
$ws.Range("H33").Value = 11500
$ws.Range("I33").Value = 5900
$ws.Range("J33").Value = 4300
$ws.Range("K33").Value = 2300
$ws.Range("L33").Value = 3300
$ws.Range("M33").Value = 3400
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("F34").Value = 0
$ws.Range("G34").Value = 0
$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = 0
$ws.Range("D35").Value = 7700
$ws.Range("E35").Value = 7300
$ws.Range("F35").Value = 7600
$ws.Range("G35").Value = 4300
$ws.Range("H35").Value = 11500
$ws.Range("I35").Value = 5900
$ws.Range("J35").Value = 4300
$ws.Range("K35").Value = 2300
$ws.Range("L35").Value = 3300
$ws.Range("M35").Value = 3400
$ws.Range("D38").Value = 43463
$ws.Range("E38").Value = 43372
$ws.Range("F38").Value = 43281
$ws.Range("G38").Value = 43190
$ws.Range("H38").Value = 43099
$ws.Range("I38").Value = 43008
$ws.Range("J38").Value = 42917
$ws.Range("K38").Value = 42826
$ws.Range("L38").Value = 42735
$ws.Range("M38").Value = 42643
$ws.Range("D41").Value = 40700
$ws.Range("E41").Value = 56200
$ws.Range("F41").Value = 16000
$ws.Range("G41").Value = 17200
$ws.Range("H41").Value = 18800
$ws.Range("I41").Value = 15600
$ws.Range("J41").Value = 19500
$ws.Range("K41").Value = 33800
$ws.Range("L41").Value = 35700
$ws.Range("M41").Value = 53100
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0
$ws.Range("F42").Value = 0
$ws.Range("G42").Value = 0
$ws.Range("H42").Value = 0
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("K42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("M42").Value = 0
$ws.Range("D43").Value = 141700
$ws.Range("E43").Value = 129500
$ws.Range("F43").Value = 120800
$ws.Range("G43").Value = 113700
$ws.Range("H43").Value = 110100
$ws.Range("I43").Value = 113000
$ws.Range("J43").Value = 104500
$ws.Range("K43").Value = 73200
$ws.Range("L43").Value = 75500
$ws.Range("M43").Value = 64100
$ws.Range("D44").Value = 0
$ws.Range("E44").Value = 0
$ws.Range("F44").Value = 0
$ws.Range("G44").Value = 0
$ws.Range("H44").Value = 0
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("M44").Value = 0
$ws.Range("D45").Value = 2600
$ws.Range("E45").Value = 2500
$ws.Range("F45").Value = 3700
$ws.Range("G45").Value = 2800
$ws.Range("H45").Value = 2600
$ws.Range("I45").Value = 2900
$ws.Range("J45").Value = 3900
$ws.Range("K45").Value = 2800
$ws.Range("L45").Value = 4000
$ws.Range("M45").Value = 2700
$ws.Range("D46").Value = 185100
$ws.Range("E46").Value = 188200
$ws.Range("F46").Value = 140600
$ws.Range("G46").Value = 133700
$ws.Range("H46").Value = 131400
$ws.Range("I46").Value = 131400
$ws.Range("J46").Value = 127800
$ws.Range("K46").Value = 109800
$ws.Range("L46").Value = 113100
$ws.Range("M46").Value = 119900
$ws.Range("D47").Value = 0
$ws.Range("E47").Value = 0
$ws.Range("F47").Value = 0
$ws.Range("G47").Value = 0
$ws.Range("H47").Value = 0
$ws.Range("I47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("K47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("M47").Value = 0
$ws.Range("D48").Value = 11700
$ws.Range("E48").Value = 10700
$ws.Range("F48").Value = 10000
$ws.Range("G48").Value = 10000
$ws.Range("H48").Value = 8700
$ws.Range("I48").Value = 8000
$ws.Range("J48").Value = 7900
$ws.Range("K48").Value = 6500
$ws.Range("L48").Value = 6700
$ws.Range("M48").Value = 3800
$ws.Range("D49").Value = 240700
$ws.Range("E49").Value = 190000
$ws.Range("F49").Value = 166000
$ws.Range("G49").Value = 168400
$ws.Range("H49").Value = 164700
$ws.Range("I49").Value = 165800
$ws.Range("J49").Value = 169600
$ws.Range("K49").Value = 97200
$ws.Range("L49").Value = 100200
$ws.Range("M49").Value = 60200
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("F50").Value = 0
$ws.Range("G50").Value = 0
$ws.Range("H50").Value = 0
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("K50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("M50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("F51").Value = 0
$ws.Range("G51").Value = 0
$ws.Range("H51").Value = 0
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("M51").Value = 0
$ws.Range("D52").Value = 2000
$ws.Range("E52").Value = 1300
$ws.Range("F52").Value = 1400
$ws.Range("G52").Value = 1000
$ws.Range("H52").Value = 1000
$ws.Range("I52").Value = 1000
$ws.Range("J52").Value = 1200
$ws.Range("K52").Value = 1400
$ws.Range("L52").Value = 1500
$ws.Range("M52").Value = 1000
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("F53").Value = 0
$ws.Range("G53").Value = 0
$ws.Range("H53").Value = 0
$ws.Range("I53").Value = 0
$ws.Range("J53").Value = 0
$ws.Range("K53").Value = 0
$ws.Range("L53").Value = 0
$ws.Range("M53").Value = 0
$ws.Range("D54").Value = 439400
$ws.Range("E54").Value = 390100
$ws.Range("F54").Value = 317900
$ws.Range("G54").Value = 313000
$ws.Range("H54").Value = 305800
$ws.Range("I54").Value = 306300
$ws.Range("J54").Value = 306500
$ws.Range("K54").Value = 214800
$ws.Range("L54").Value = 221500
$ws.Range("M54").Value = 184800
$ws.Range("D57").Value = 22700
$ws.Range("E57").Value = 17900
$ws.Range("F57").Value = 18900
$ws.Range("G57").Value = 17100
$ws.Range("H57").Value = 18700
$ws.Range("I57").Value = 17400
$ws.Range("J57").Value = 19100
$ws.Range("K57").Value = 9700
$ws.Range("L57").Value = 13700
$ws.Range("M57").Value = 11800
$ws.Range("D58").Value = 17100
$ws.Range("E58").Value = 14200
$ws.Range("F58").Value = 12700
$ws.Range("G58").Value = 12400
$ws.Range("H58").Value = 11100
$ws.Range("I58").Value = 10800
$ws.Range("J58").Value = 11000
$ws.Range("K58").Value = 11800
$ws.Range("L58").Value = 11300
$ws.Range("M58").Value = 6600
$ws.Range("D59").Value = 33100
$ws.Range("E59").Value = 22400
$ws.Range("F59").Value = 25200
$ws.Range("G59").Value = 24800
$ws.Range("H59").Value = 26600
$ws.Range("I59").Value = 23400
$ws.Range("J59").Value = 24300
$ws.Range("K59").Value = 14400
$ws.Range("L59").Value = 18600
$ws.Range("M59").Value = 13900
$ws.Range("D60").Value = 73000
$ws.Range("E60").Value = 54600
$ws.Range("F60").Value = 56800
$ws.Range("G60").Value = 54300
$ws.Range("H60").Value = 56400
$ws.Range("I60").Value = 51700
$ws.Range("J60").Value = 54500
$ws.Range("K60").Value = 35900
$ws.Range("L60").Value = 43600
$ws.Range("M60").Value = 32400
$ws.Range("D61").Value = 29800
$ws.Range("E61").Value = 19000
$ws.Range("F61").Value = 50300
$ws.Range("G61").Value = 59000
$ws.Range("H61").Value = 57400
$ws.Range("I61").Value = 67200
$ws.Range("J61").Value = 71800
$ws.Range("K61").Value = 21300
$ws.Range("L61").Value = 21600
$ws.Range("M61").Value = 10000
$ws.Range("D62").Value = 19100
$ws.Range("E62").Value = 13200
$ws.Range("F62").Value = 13100
$ws.Range("G62").Value = 12700
$ws.Range("H62").Value = 11800
$ws.Range("I62").Value = 22200
$ws.Range("J62").Value = 22200
$ws.Range("K62").Value = 6500
$ws.Range("L62").Value = 8100
$ws.Range("M62").Value = 2200
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("F63").Value = 0
$ws.Range("G63").Value = 0
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("F64").Value = 0
$ws.Range("G64").Value = 0
$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("M64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("F65").Value = 0
$ws.Range("G65").Value = 0
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = 0
$ws.Range("D66").Value = 121900
$ws.Range("E66").Value = 86800
$ws.Range("F66").Value = 120100
$ws.Range("G66").Value = 126000
$ws.Range("H66").Value = 125700
$ws.Range("I66").Value = 141000
$ws.Range("J66").Value = 148500
$ws.Range("K66").Value = 63600
$ws.Range("L66").Value = 73300
$ws.Range("M66").Value = 44500
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("F68").Value = 0
$ws.Range("G68").Value = 0
$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("F69").Value = 0
$ws.Range("G69").Value = 0
$ws.Range("H69").Value = 0
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("M69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("F70").Value = 0
$ws.Range("G70").Value = 0
$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("F71").Value = 0
$ws.Range("G71").Value = 0
$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = 0
$ws.Range("D72").Value = 80900
$ws.Range("E72").Value = 73200
$ws.Range("F72").Value = 65900
$ws.Range("G72").Value = 58300
$ws.Range("H72").Value = 54000
$ws.Range("I72").Value = 42500
$ws.Range("J72").Value = 36600
$ws.Range("K72").Value = 32300
$ws.Range("L72").Value = 30000
$ws.Range("M72").Value = 26700
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("F73").Value = 0
$ws.Range("G73").Value = 0
$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("F74").Value = 0
$ws.Range("G74").Value = 0
$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("F75").Value = 0
$ws.Range("G75").Value = 0
$ws.Range("H75").Value = 0
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("M75").Value = 0
$ws.Range("D76").Value = 317500
$ws.Range("E76").Value = 303300
$ws.Range("F76").Value = 197800
$ws.Range("G76").Value = 187000
$ws.Range("H76").Value = 180100
$ws.Range("I76").Value = 165200
$ws.Range("J76").Value = 158100
$ws.Range("K76").Value = 151200
$ws.Range("L76").Value = 148200
$ws.Range("M76").Value = 140300
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("F77").Value = 0
$ws.Range("G77").Value = 0
$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = 0
$ws.Range("D80").Value = 43463
$ws.Range("E80").Value = 43372
$ws.Range("F80").Value = 43281
$ws.Range("G80").Value = 43190
$ws.Range("H80").Value = 43099
$ws.Range("I80").Value = 43008
$ws.Range("J80").Value = 42917
$ws.Range("K80").Value = 42826
$ws.Range("L80").Value = 42735
$ws.Range("M80").Value = 42643
$ws.Range("D81").Value = 7700
$ws.Range("E81").Value = 7300
$ws.Range("F81").Value = 7600
$ws.Range("G81").Value = 4300
$ws.Range("H81").Value = 11500
$ws.Range("I81").Value = 5900
$ws.Range("J81").Value = 4300
$ws.Range("K81").Value = 2300
$ws.Range("L81").Value = 3300
$ws.Range("M81").Value = 3400
$ws.Range("D83").Value = 5700
$ws.Range("E83").Value = 4100
$ws.Range("F83").Value = 3800
$ws.Range("G83").Value = 3800
$ws.Range("H83").Value = 3600
$ws.Range("I83").Value = 3800
$ws.Range("J83").Value = 3300
$ws.Range("K83").Value = 2500
$ws.Range("L83").Value = 1900
$ws.Range("M83").Value = 1600
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("F84").Value = 0
$ws.Range("G84").Value = 0
$ws.Range("H84").Value = 0
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("F85").Value = 0
$ws.Range("G85").Value = 0
$ws.Range("H85").Value = 0
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("F86").Value = 0
$ws.Range("G86").Value = 0
$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("F87").Value = 0
$ws.Range("G87").Value = 0
$ws.Range("H87").Value = 0
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("K87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("M87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("F88").Value = 0
$ws.Range("G88").Value = 0
$ws.Range("H88").Value = 0
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("M88").Value = 0
$ws.Range("D89").Value = 17400
$ws.Range("E89").Value = 6700
$ws.Range("F89").Value = 7200
$ws.Range("G89").Value = 3600
$ws.Range("H89").Value = 12000
$ws.Range("I89").Value = 1200
$ws.Range("J89").Value = 4300
$ws.Range("K89").Value = 200
$ws.Range("L89").Value = 4500
$ws.Range("M89").Value = 4000
$ws.Range("D91").Value = -600
$ws.Range("E91").Value = -100
$ws.Range("F91").Value = -300
$ws.Range("G91").Value = -1100
$ws.Range("H91").Value = -600
$ws.Range("I91").Value = -600
$ws.Range("J91").Value = -600
$ws.Range("K91").Value = -400
$ws.Range("L91").Value = -400
$ws.Range("M91").Value = -100
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("F92").Value = 0
$ws.Range("G92").Value = 0
$ws.Range("H92").Value = 0
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("M92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("F93").Value = 0
$ws.Range("G93").Value = 0
$ws.Range("H93").Value = 0
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = 0
$ws.Range("D94").Value = -30300
$ws.Range("E94").Value = -25100
$ws.Range("F94").Value = -500
$ws.Range("G94").Value = -4400
$ws.Range("H94").Value = -1000
$ws.Range("I94").Value = 200
$ws.Range("J94").Value = -61600
$ws.Range("K94").Value = -400
$ws.Range("L94").Value = -21800
$ws.Range("M94").Value = -400
$ws.Range("D96").Value = 0
$ws.Range("E96").Value = 0
$ws.Range("F96").Value = 0
$ws.Range("G96").Value = 0
$ws.Range("H96").Value = 0
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("M96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("F97").Value = 0
$ws.Range("G97").Value = 0
$ws.Range("H97").Value = 0
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("F98").Value = 0
$ws.Range("G98").Value = 0
$ws.Range("H98").Value = 0
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("F99").Value = 0
$ws.Range("G99").Value = 0
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = 0
$ws.Range("D100").Value = -2600
$ws.Range("E100").Value = 58500
$ws.Range("F100").Value = -7900
$ws.Range("G100").Value = -700
$ws.Range("H100").Value = -7700
$ws.Range("I100").Value = -5300
$ws.Range("J100").Value = 43000
$ws.Range("K100").Value = -1600
$ws.Range("L100").Value = -100
$ws.Range("M100").Value = -1500
$ws.Range("D101").Value = 0
$ws.Range("E101").Value = 0
$ws.Range("F101").Value = 0
$ws.Range("G101").Value = 0
$ws.Range("H101").Value = 0
$ws.Range("I101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("K101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("M101").Value = 0
$ws.Range("D102").Value = -15500
$ws.Range("E102").Value = 40200
$ws.Range("F102").Value = -1200
$ws.Range("G102").Value = -1500
$ws.Range("H102").Value = 3200
$ws.Range("I102").Value = -3900
$ws.Range("J102").Value = -14400
$ws.Range("K102").Value = -1800
$ws.Range("L102").Value = -17400
$ws.Range("M102").Value = 2000
